# "Generate Report for Handback" - localization-status.xlsx update
#
# Applies:
#  - "In Translation" -> "Handed back: in sync with en-US" everywhere (Overview + both language sheets)
#  - zh-cn sheet: Latest Handback DateTime (K2/K3) "0001-01-01 00:00:00" -> "2016-09-06 06:28:33"
#  - de-de sheet: Latest Handback DateTime (K2/K3) -> new timestamp "2016-09-06 06:28:41"
#  - zh-cn / de-de sheets: populate "Latest Target File" (I) with a hyperlink to the handback
#    markdown (same target as column A's link) and "Latest Handback File" (J) with the xlf name
#  - widen columns that now hold the longer handback text / links

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (appears on Overview!E2:F3 and on both language sheets' Status column C2:C3)
# ---------------------------------------------------------------------------
$overview.Cells.Replace("In Translation", "Handed back: in sync with en-US")
$zhcn.Cells.Replace("In Translation", "Handed back: in sync with en-US")
$dede.Cells.Replace("In Translation", "Handed back: in sync with en-US")

# ---------------------------------------------------------------------------
# 2. Handback target URLs - reuse the same external links already used by
#    column A ("…md" files on GitHub) for the new "Latest Target File" links.
# ---------------------------------------------------------------------------
$url638 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/696b912a0949c09cb9572acb68e7e510d8c609f7/e2e/638b9d65-6654-429d-aa80-88cda2ef2088.md"
$url_e6 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/696b912a0949c09cb9572acb68e7e510d8c609f7/e2e/e6be2e76-0ffa-4f65-9675-710600040798.md"

# ---------------------------------------------------------------------------
# 3. zh-cn sheet ("Latest Target File" + "Latest Handback File" + DateTime)
# ---------------------------------------------------------------------------
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $url638, "", "", "638b9d65-6654-429d-aa80-88cda2ef2088.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $url638, "", "", "638b9d65-6654-429d-aa80-88cda2ef2088.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $url_e6, "", "", "e6be2e76-0ffa-4f65-9675-710600040798.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $url_e6, "", "", "e6be2e76-0ffa-4f65-9675-710600040798.md")

$zhcn.Range("J2").Value = "638b9d65-6654-429d-aa80-88cda2ef2088.b28b39e35365431b33c30995d9b157d12505d088.zh-cn.xlf"
$zhcn.Range("J3").Value = "e6be2e76-0ffa-4f65-9675-710600040798.a67a5f49fc677cfb818bb0e22113c87300a2ab13.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-09-06 06:28:33"
$zhcn.Range("K3").Value = "2016-09-06 06:28:33"

# ---------------------------------------------------------------------------
# 4. de-de sheet (same shape, different handback timestamp)
# ---------------------------------------------------------------------------
$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $url638, "", "", "638b9d65-6654-429d-aa80-88cda2ef2088.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $url638, "", "", "638b9d65-6654-429d-aa80-88cda2ef2088.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $url_e6, "", "", "e6be2e76-0ffa-4f65-9675-710600040798.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $url_e6, "", "", "e6be2e76-0ffa-4f65-9675-710600040798.md")

$dede.Range("J2").Value = "638b9d65-6654-429d-aa80-88cda2ef2088.b28b39e35365431b33c30995d9b157d12505d088.de-de.xlf"
$dede.Range("J3").Value = "e6be2e76-0ffa-4f65-9675-710600040798.a67a5f49fc677cfb818bb0e22113c87300a2ab13.de-de.xlf"

$dede.Range("K2").Value = "2016-09-06 06:28:41"
$dede.Range("K3").Value = "2016-09-06 06:28:41"

# ---------------------------------------------------------------------------
# 5. Column widths - widened to fit the new handback text / hyperlink text.
#    ColumnWidth is expressed in characters and Excel stores/rounds it to the
#    nearest pixel, so we feed it the (width - 5/6) pre-image of each target
#    "stored" width (29.9777... rounds to the 30-character pixel bucket, 40
#    lands exactly on a pixel bucket already).
# ---------------------------------------------------------------------------
$wideStatus = 29.166666666666668   # -> stored width ~29.98 (Status/handback columns)
$wideLink   = 39.166666666666664   # -> stored width 40     (Target/Handback file columns)

$overview.Columns.Item(5).ColumnWidth = $wideStatus   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = $wideStatus   # F: de-de status

$zhcn.Columns.Item(3).ColumnWidth  = $wideStatus       # C: Status
$zhcn.Columns.Item(9).ColumnWidth  = $wideLink         # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = $wideLink         # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth  = $wideStatus       # C: Status
$dede.Columns.Item(9).ColumnWidth  = $wideLink         # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = $wideLink         # J: Latest Handback File
